$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 'Tứ'
$ws.Cells.Item(2, 3).Value = 'Lê Nhã'
$ws.Cells.Item(2, 4).Value = 'Quảng Nam'
$ws.Cells.Item(3, 2).Value = 'Lan'
$ws.Cells.Item(3, 3).Value = 'Trần Huỳnh'
$ws.Cells.Item(3, 4).Value = 'Đà Nẵng'
$ws.Cells.Item(4, 2).Value = 'Huệ'
$ws.Cells.Item(4, 3).Value = 'Nguyễn Công'
$ws.Cells.Item(4, 4).Value = 'Hồ Chí Minh'
$ws.Cells.Item(5, 2).Value = 'Hoa'
$ws.Cells.Item(5, 3).Value = 'Huỳnh Như'
$ws.Cells.Item(5, 4).Value = 'Đồng Nai'
$ws.Cells.Item(6, 2).Value = 'Lộc'
$ws.Cells.Item(6, 3).Value = 'Hoàng Nhã'
$ws.Cells.Item(6, 4).Value = 'Đà Nẵng'
$ws.Cells.Item(7, 2).Value = 'Huyền'
$ws.Cells.Item(7, 3).Value = 'Nguyễn Gia'
$ws.Cells.Item(7, 4).Value = 'Đà Nẵng'
$ws.Cells.Item(8, 2).Value = 'Hoa'
$ws.Cells.Item(8, 3).Value = 'Phạm Ngọc'
$ws.Cells.Item(8, 4).Value = 'Hải Dương'
$ws.Cells.Item(9, 2).Value = 'Hổ'
$ws.Cells.Item(9, 3).Value = 'Hoàng Văn'
$ws.Cells.Item(9, 4).Value = 'Hồ Chí Minh'
$ws.Cells.Item(10, 2).Value = 'Tứ'
$ws.Cells.Item(10, 3).Value = 'Phạm Văn'
$ws.Cells.Item(10, 4).Value = 'Hà Nội'
$ws.Cells.Item(11, 2).Value = 'Lan'
$ws.Cells.Item(11, 3).Value = 'Nguyễn Gia'
$ws.Cells.Item(11, 4).Value = 'Hải Dương'
$ws.Cells.Item(12, 2).Value = 'Nhị'
$ws.Cells.Item(12, 3).Value = 'Đào Như'
$ws.Cells.Item(12, 4).Value = 'Đà Nẵng'
$ws.Cells.Item(13, 2).Value = 'Huệ'
$ws.Cells.Item(13, 3).Value = 'Phạm Như'
$ws.Cells.Item(13, 4).Value = 'Quảng Ninh'
$ws.Cells.Item(14, 2).Value = 'Lộc'
$ws.Cells.Item(14, 3).Value = 'Lê Công'
$ws.Cells.Item(14, 4).Value = 'Quảng Nam'
$ws.Cells.Item(15, 2).Value = 'Huệ'
$ws.Cells.Item(15, 3).Value = 'Lê Văn'
$ws.Cells.Item(15, 4).Value = 'Vũng Tàu'
$ws.Cells.Item(16, 2).Value = 'Tam'
$ws.Cells.Item(16, 3).Value = 'Trần Văn'
$ws.Cells.Item(16, 4).Value = 'Vũng Tàu'
$ws.Cells.Item(17, 2).Value = 'Lan'
$ws.Cells.Item(17, 3).Value = 'Đào Nhã'
$ws.Cells.Item(17, 4).Value = 'Hải Dương'
$ws.Cells.Item(18, 2).Value = 'Ngũ'
$ws.Cells.Item(18, 3).Value = 'Đào Văn'
$ws.Cells.Item(18, 4).Value = 'Hải Dương'
$ws.Cells.Item(19, 2).Value = 'Lộc'
$ws.Cells.Item(19, 3).Value = 'Trần Công'
$ws.Cells.Item(19, 4).Value = 'Quảng Ninh'
$ws.Cells.Item(20, 2).Value = 'Lộc'
$ws.Cells.Item(20, 3).Value = 'Nguyễn Nhã'
$ws.Cells.Item(20, 4).Value = 'Huế'
$ws.Cells.Item(21, 2).Value = 'Nhị'
$ws.Cells.Item(21, 3).Value = 'Lê Huỳnh'
$ws.Cells.Item(21, 4).Value = 'Hà Nội'
$ws.Cells.Item(22, 2).Value = 'Ngũ'
$ws.Cells.Item(22, 3).Value = 'Trần Ngọc'
$ws.Cells.Item(22, 4).Value = 'Hải Dương'
$ws.Cells.Item(23, 2).Value = 'Ngũ'
$ws.Cells.Item(23, 3).Value = 'Lê Văn'
$ws.Cells.Item(23, 4).Value = 'Vũng Tàu'
$ws.Cells.Item(24, 2).Value = 'Nhị'
$ws.Cells.Item(24, 3).Value = 'Lê Văn'
$ws.Cells.Item(24, 4).Value = 'Đà Nẵng'
$ws.Cells.Item(25, 2).Value = 'Ngũ'
$ws.Cells.Item(25, 3).Value = 'Nguyễn Huỳnh'
$ws.Cells.Item(25, 4).Value = 'Vũng Tàu'
$ws.Cells.Item(26, 2).Value = 'Nhị'
$ws.Cells.Item(26, 3).Value = 'Huỳnh Huỳnh'
$ws.Cells.Item(26, 4).Value = 'Đà Nẵng'
$ws.Cells.Item(27, 2).Value = 'Nhị'
$ws.Cells.Item(27, 3).Value = 'Huỳnh Nhã'
$ws.Cells.Item(27, 4).Value = 'Hải Dương'
$ws.Cells.Item(28, 2).Value = 'Huyền'
$ws.Cells.Item(28, 3).Value = 'Đào Huỳnh'
$ws.Cells.Item(28, 4).Value = 'Quảng Nam'
$ws.Cells.Item(29, 2).Value = 'Lộc'
$ws.Cells.Item(29, 3).Value = 'Huỳnh Gia'
$ws.Cells.Item(29, 4).Value = 'Hà Nội'
$ws.Cells.Item(30, 2).Value = 'Tam'
$ws.Cells.Item(30, 3).Value = 'Lê Ngọc'
$ws.Cells.Item(30, 4).Value = 'Huế'
$ws.Cells.Item(31, 2).Value = 'Lan'
$ws.Cells.Item(31, 3).Value = 'Phạm Ngọc'
$ws.Cells.Item(31, 4).Value = 'Vũng Tàu'
$ws.Cells.Item(32, 2).Value = 'Nhị'
$ws.Cells.Item(32, 3).Value = 'Hoàng Như'
$ws.Cells.Item(32, 4).Value = 'Đà Nẵng'
$ws.Cells.Item(33, 2).Value = 'Tứ'
$ws.Cells.Item(33, 3).Value = 'Nguyễn Huỳnh'
$ws.Cells.Item(33, 4).Value = 'Đà Nẵng'
$ws.Cells.Item(34, 2).Value = 'Ngũ'
$ws.Cells.Item(34, 3).Value = 'Đào Nhã'
$ws.Cells.Item(34, 4).Value = 'Hà Nội'
$ws.Cells.Item(35, 2).Value = 'Nhị'
$ws.Cells.Item(35, 3).Value = 'Lê Công'
$ws.Cells.Item(35, 4).Value = 'Quảng Ninh'
$ws.Cells.Item(36, 2).Value = 'Tứ'
$ws.Cells.Item(36, 3).Value = 'Lê Nhã'
$ws.Cells.Item(36, 4).Value = 'Đồng Nai'
$ws.Cells.Item(37, 2).Value = 'Tam'
$ws.Cells.Item(37, 3).Value = 'Huỳnh Văn'
$ws.Cells.Item(37, 4).Value = 'Vũng Tàu'
$ws.Cells.Item(38, 2).Value = 'Ngũ'
$ws.Cells.Item(38, 3).Value = 'Lê Công'
$ws.Cells.Item(38, 4).Value = 'Quảng Ninh'
$ws.Cells.Item(39, 2).Value = 'Cúc'
$ws.Cells.Item(39, 3).Value = 'Huỳnh Huỳnh'
$ws.Cells.Item(39, 4).Value = 'Hải Dương'
$ws.Cells.Item(40, 2).Value = 'Huệ'
$ws.Cells.Item(40, 3).Value = 'Đào Nhã'
$ws.Cells.Item(40, 4).Value = 'Hồ Chí Minh'
$ws.Cells.Item(41, 2).Value = 'Huyền'
$ws.Cells.Item(41, 3).Value = 'Trần Huỳnh'
$ws.Cells.Item(41, 4).Value = 'Đồng Nai'
$ws.Cells.Item(42, 2).Value = 'Tam'
$ws.Cells.Item(42, 3).Value = 'Phạm Gia'
$ws.Cells.Item(42, 4).Value = 'Hải Dương'
$ws.Cells.Item(43, 2).Value = 'Huyền'
$ws.Cells.Item(43, 3).Value = 'Nguyễn Nhã'
$ws.Cells.Item(43, 4).Value = 'Đồng Nai'
$ws.Cells.Item(44, 2).Value = 'Tam'
$ws.Cells.Item(44, 3).Value = 'Hoàng Công'
$ws.Cells.Item(44, 4).Value = 'Huế'
$ws.Cells.Item(45, 2).Value = 'Hoa'
$ws.Cells.Item(45, 3).Value = 'Lê Nhã'
$ws.Cells.Item(45, 4).Value = 'Quảng Nam'
$ws.Cells.Item(46, 2).Value = 'Tứ'
$ws.Cells.Item(46, 3).Value = 'Hoàng Gia'
$ws.Cells.Item(46, 4).Value = 'Quảng Ninh'
$ws.Cells.Item(47, 2).Value = 'Tam'
$ws.Cells.Item(47, 3).Value = 'Nguyễn Ngọc'
$ws.Cells.Item(47, 4).Value = 'Hải Dương'
$ws.Cells.Item(48, 2).Value = 'Huyền'
$ws.Cells.Item(48, 3).Value = 'Đào Ngọc'
$ws.Cells.Item(48, 4).Value = 'Hải Dương'
$ws.Cells.Item(49, 2).Value = 'Tứ'
$ws.Cells.Item(49, 3).Value = 'Lê Gia'
$ws.Cells.Item(49, 4).Value = 'Quảng Nam'
$ws.Cells.Item(50, 2).Value = 'Ngũ'
$ws.Cells.Item(50, 3).Value = 'Phạm Công'
$ws.Cells.Item(50, 4).Value = 'Hồ Chí Minh'
$ws.Cells.Item(51, 2).Value = 'Tam'
$ws.Cells.Item(51, 3).Value = 'Trần Huỳnh'
$ws.Cells.Item(51, 4).Value = 'Hà Nội'

$ws.Columns.Item(3).ColumnWidth = 16.0
$ws.Columns.Item(4).ColumnWidth = 13.2
